$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1.
# Copy formatting from the existing header cell (AC1) first so the new
# headers pick up the same bold/centered/bordered header style, then
# overwrite the text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season-record columns for every data row (2-51): every row
# in this sheet corresponds to the same team/season, so the record is
# identical (72 wins, 90 losses, 0 ties) on each row.
$ws.Range("AD2:AD51").Value = 72
$ws.Range("AE2:AE51").Value = 90
$ws.Range("AF2:AF51").Value = 0
